$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# The "merge_id" info that used to live in a comment on Report!A7 is now
# appended directly to the cell's own text, so the template no longer
# needs a separate comment to carry it.
$ws1.Range("A7").Value = "<#table.NAME_GROUP><#merge.<#table.NAME_GROUP>>"
$ws1.Range("A7").Comment.Delete()

# Column A on the Report sheet got a bit wider to fit the longer text.
$ws1.Columns.Item(1).ColumnWidth = 26.5

# Rows that used to carry an explicit (slightly taller) row height now
# just use the sheet's default row height.
foreach ($rowNum in @(1,5,6,7,8,10,12)) {
    $ws1.Rows.Item($rowNum).EntireRow.AutoFit()
}

# Update the remembered selection on each sheet (select sheet1's range
# last so the Report sheet stays the active tab).
[void]$ws2.Range("A4").Select()
[void]$ws1.Range("B9").Select()
